$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "DAC101C081" worksheet after the last sheet (TMF8801)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "DAC101C081"

# Column A width (matches the other register sheets' "Name" column)
$ws.Columns.Item(1).ColumnWidth = 20.09

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Hex Address"
$ws.Range("C1").Value = "Default Value"
$ws.Range("D1").Value = "Bit Width"
$ws.Range("E1").Value = "Bit Index (High)"
$ws.Range("F1").Value = "Bit Index (Low)"

# Row 2 - PD field
$ws.Range("A2").Value = "PD"
$ws.Range("B2").Value = "None"
$ws.Range("C2").Value = "0x0000"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 12

# Row 3 - DATA field
$ws.Range("A3").Value = "DATA"
$ws.Range("B3").Value = "None"
$ws.Range("C3").Value = "0x0000"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = 2

# The "Default Value" column uses a black-on-no-fill variant of the built in
# "Bad" cell style (normal text color, no red fill).
$ws.Range("C2:C3").Font.ThemeColor = 1
$ws.Range("C2:C3").Interior.ColorIndex = 0

$ws.Range("B3").Select()

# ---------------------------------------------------------------------------
# 2. DAC53401 sheet - the active selection becomes the header range
# ---------------------------------------------------------------------------
$dac53401 = $wb.Worksheets.Item("DAC53401")
$dac53401.Range("A1:F3").Select()

# Re-activate the newly added sheet so it ends up as the active tab, matching
# the target workbook state.
$ws.Activate()
$ws.Range("B3").Select()
